$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two now-unwanted product rows (rows 6 and 7), bottom-up so
# row indices of remaining rows don't shift while we work.
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()

# Update the Images column (F) for the remaining data rows: instead of the
# old comma separated "cintamani275_7, ..." placeholders, just record the
# actual image filename(s) that should exist on disk before import.
$ws.Range("F3").Value = "002.jpg"
$ws.Range("F2").Value = "001.JPG"
$ws.Range("F4").Value = "001.JPG"
$ws.Range("F5").Value = "001.JPG"

# Format column F (Images) as Text so filenames like "001.JPG" aren't
# reinterpreted, and move the active selection to match.
$ws.Range("F1:F5").NumberFormat = "@"

$ws.Range("F5").Select()
